# Fall 2021 LLO8180 due dates - schedule update
#
# The due-date schedule slipped by roughly one week starting with the
# Week 8 quiz (and there's an extra week's gap right before the final
# week, presumably for Thanksgiving break). The final "Last class" date
# also got pinned down (it was "TBD") and its note text was clarified.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Due dates move later (serial date numbers, one week = 7 days) ---
$ws.Range("B14").Value = 44475   # Q5  due date:  9/29/21 -> 10/6/21
$ws.Range("B16").Value = 44482   # Q6  due date: 10/6/21  -> 10/13/21
$ws.Range("B19").Value = 44489   # Q7  due date: 10/13/21 -> 10/20/21
$ws.Range("B21").Value = 44496   # Q8  due date: 10/20/21 -> 10/27/21
$ws.Range("B24").Value = 44503   # Q9  due date: 10/27/21 -> 11/3/21
$ws.Range("B26").Value = 44510   # Q10 due date: 11/3/21  -> 11/10/21
$ws.Range("B29").Value = 44517   # Q11 due date: 11/10/21 -> 11/17/21
$ws.Range("B31").Value = 44531   # Q12 due date: 11/17/21 -> 12/1/21  (two week gap)
$ws.Range("B34").Value = 44538   # Q13 due date: 12/1/21  -> 12/8/21

# B36 (A7 / Q14 row) previously held the literal text "TBD" for the
# final class date - it's now a real date.
$ws.Range("B36").Value = 44545   # 12/15/21

# C34 note text changes from "Last class!" to "Last class meeting!"
$ws.Range("C34").Value = "Last class meeting!"

# C36:C37 ("due by midnight" note, merged cell) picks up the same
# bold/italic centered style used elsewhere, but now with the thin
# left border that separates it from column B.
$ws.Range("C36:C37").HorizontalAlignment = -4108   # xlCenter
$ws.Range("C36:C37").VerticalAlignment = -4108     # xlCenter
$ws.Range("C36:C37").Borders.Item(7).LineStyle = 1 # xlContinuous (left border)

# Column C widened to fit the longer "Last class meeting!" text.
$ws.Columns.Item(3).ColumnWidth = 17.6

# Cursor/selection left on D33 (matches the saved workbook state).
[void]$ws.Range("D33").Select()
